# "Simplified data processing and added YPD data"
# Adds a new "Dataset id" column (H) to the phenotype-mapping sheet,
# tagging each existing row with the Yeastphenome dataset id it came from.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header cell, styled like the rest of row 1 (bold).
$ws.Range("H1").Value = "Dataset id"
$ws.Range("H1").Font.Bold = $true

# Dataset id values for the six data rows, one per original row (2-7).
$datasetIds = @(725, 725, 726, 727, 729, 728)
for ($i = 0; $i -lt $datasetIds.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $datasetIds[$i]
}

# Move the selection to the cell below the newly-entered data, as Excel
# would leave it after typing the last value and pressing Enter.
$ws.Range("H8").Select() | Out-Null
